# Update "想去人数" (want-to-go count) values in both the "展览" and
# "全部类型" worksheets, which carry duplicate data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F5").Value = 2219
    $ws.Range("F7").Value = 1363
    $ws.Range("F9").Value = 144
    $ws.Range("F11").Value = 325
}
